# Update the cached "datetimeFigureOut" date field text that appears on the
# slide master, every slide layout, and the notes master (Insert > Header &
# Footer "Date and time" placeholder) from 7/29/2017 to 10/19/2017.

$p = $ppt.ActivePresentation

$oldDate = "7/29/2017"
$newDate = "10/19/2017"

function Update-DatePlaceholder($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $ph = $sh.PlaceholderFormat
            if ($ph -ne $null -and $ph.Type -eq 16) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master.
Update-DatePlaceholder $p.SlideMaster

# Every slide layout attached to the slide master.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L)
}

# Notes master.
Update-DatePlaceholder $p.NotesMaster
